$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-26 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2026-01-27 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("40÷7=5, 5", $true, $true, $false, $false, $false, $true, 1, $false, "43÷9=4, 7", 2) | Out-Null
$d.Content.Find.Execute("83÷6=13, 5", $true, $true, $false, $false, $false, $true, 1, $false, "50÷3=16, 2", 2) | Out-Null
$d.Content.Find.Execute("87÷5=17, 2", $true, $true, $false, $false, $false, $true, 1, $false, "25÷3=8, 1", 2) | Out-Null
$d.Content.Find.Execute("21÷2=10, 1", $true, $true, $false, $false, $false, $true, 1, $false, "97÷4=24, 1", 2) | Out-Null
$d.Content.Find.Execute("10÷3=3, 1", $true, $true, $false, $false, $false, $true, 1, $false, "68÷5=13, 3", 2) | Out-Null
$d.Content.Find.Execute("48÷8=6, 0", $true, $true, $false, $false, $false, $true, 1, $false, "80÷6=13, 2", 2) | Out-Null
$d.Content.Find.Execute("55÷8=6, 7", $true, $true, $false, $false, $false, $true, 1, $false, "34÷9=3, 7", 2) | Out-Null
$d.Content.Find.Execute("65÷7=9, 2", $true, $true, $false, $false, $false, $true, 1, $false, "33÷7=4, 5", 2) | Out-Null
$d.Content.Find.Execute("79÷5=15, 4", $true, $true, $false, $false, $false, $true, 1, $false, "24÷2=12, 0", 2) | Out-Null
$d.Content.Find.Execute("86÷5=17, 1", $true, $true, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 2) | Out-Null
$d.Content.Find.Execute("65÷6=10, 5", $true, $true, $false, $false, $false, $true, 1, $false, "67÷3=22, 1", 2) | Out-Null
$d.Content.Find.Execute("81÷2=40, 1", $true, $true, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2) | Out-Null
$d.Content.Find.Execute("19÷5=3, 4", $true, $true, $false, $false, $false, $true, 1, $false, "32÷5=6, 2", 2) | Out-Null
$d.Content.Find.Execute("77÷6=12, 5", $true, $true, $false, $false, $false, $true, 1, $false, "12÷7=1, 5", 2) | Out-Null
$d.Content.Find.Execute("22÷7=3, 1", $true, $true, $false, $false, $false, $true, 1, $false, "32÷5=6, 2", 2) | Out-Null
$d.Content.Find.Execute("64÷6=10, 4", $true, $true, $false, $false, $false, $true, 1, $false, "76÷3=25, 1", 2) | Out-Null
$d.Content.Find.Execute("63÷2=31, 1", $true, $true, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 2) | Out-Null
$d.Content.Find.Execute("91÷5=18, 1", $true, $true, $false, $false, $false, $true, 1, $false, "20÷3=6, 2", 2) | Out-Null
$d.Content.Find.Execute("18÷3=6, 0", $true, $true, $false, $false, $false, $true, 1, $false, "57÷5=11, 2", 2) | Out-Null
$d.Content.Find.Execute("88÷2=44, 0", $true, $true, $false, $false, $false, $true, 1, $false, "27÷8=3, 3", 2) | Out-Null
$d.Content.Find.Execute("33÷2=16, 1", $true, $true, $false, $false, $false, $true, 1, $false, "24÷3=8, 0", 2) | Out-Null
$d.Content.Find.Execute("58÷8=7, 2", $true, $true, $false, $false, $false, $true, 1, $false, "88÷8=11, 0", 2) | Out-Null
$d.Content.Find.Execute("78÷4=19, 2", $true, $true, $false, $false, $false, $true, 1, $false, "69÷5=13, 4", 2) | Out-Null
$d.Content.Find.Execute("63÷8=7, 7", $true, $true, $false, $false, $false, $true, 1, $false, "24÷2=12, 0", 2) | Out-Null
$d.Content.Find.Execute("27÷4=6, 3", $true, $true, $false, $false, $false, $true, 1, $false, "68÷4=17, 0", 2) | Out-Null
